$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.194.55'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '2.935.03'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '356.50'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').Value = '109.67'
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('E7').Value = '  +1.79%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '0.627'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').Value = '38.97'
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').Value = '0.0874'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').Value = '19.52'
$ws.Range('E13').Value = '  -1.67%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.406.11'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '7.76'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '2.939.93'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').Value = '0.979'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').Value = '52.229.79'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '3.54'
$ws.Range('E19').Value = '  +7.20%  '
$ws.Range('D20').Value = '7.58'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').Value = '13.88'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = '0.0₃0980'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').Value = '70.44'
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').Value = '270.48'
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').Value = '2.80'
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('D26').Value = '0.179'
$ws.Range('E26').Value = '  +3.04%  '
$ws.Range('D27').Value = '7.90'
$ws.Range('E27').Value = '  +19.01%  '
$ws.Range('D28').Value = '26.95'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  +9.12%  '
$ws.Range('D31').Value = '10.46'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '2.28'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = '37.51'
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').Value = '6.21'
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').Value = '51.98'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('D36').Value = '0.0444'
$ws.Range('E36').Value = '  -2.11%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '3.19'
$ws.Range('E38').Value = '  -5.46%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '18.21'
$ws.Range('E39').Value = '  -4.11%  '
$ws.Range('D40').Value = '1.99'
$ws.Range('E40').Value = '  -3.73%  '
$ws.Range('D41').Value = '2.74'
$ws.Range('E41').Value = '  -2.52%  '
$ws.Range('E42').Value = '  +2.09%  '
$ws.Range('D43').Value = '23.02'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').Value = '119.74'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('D46').Value = '3.46'
$ws.Range('E46').Value = '  -2.03%  '
$ws.Range('E47').Value = '  -5.41%  '
$ws.Range('D48').Value = '2.133.62'
$ws.Range('E48').Value = '  -2.54%  '
$ws.Range('D49').Value = '0.250'
$ws.Range('E49').Value = '  -5.60%  '
$ws.Range('D50').Value = '0.0353'
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').Value = '62.31'
$ws.Range('E51').Value = '  +2.29%  '
